$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text: "11/02/2025" -> "13/02/2025" on the slide master
#    and on every slide layout (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------------
$newDate = "13/02/2025"

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shp = $master.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $isDate = $false
        try { $isDate = ($shp.PlaceholderFormat.Type -eq 16) } catch { $isDate = $false }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lyt = $layouts.Item($i)
    for ($j = 1; $j -le $lyt.Shapes.Count; $j++) {
        $shp = $lyt.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try { $isDate = ($shp.PlaceholderFormat.Type -eq 16) } catch { $isDate = $false }
            if ($isDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 3 table: rewording "Good"/"Medium" -> "High"/"Moderate"
#    (first version of article and SM)
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tbl = $s3.Shapes.Item(1).Table

# Row 2 (Medical records): Quality "Good " -> "High ", Quantity "Good" -> "High"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Characters(1, 4).Text = "High"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Characters(1, 4).Text = "High"

# Row 3 (CRF / self-questionnaires): "Medium: " -> "Moderate: " in both columns
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Characters(1, 8).Text = "Moderate: "
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Characters(1, 8).Text = "Moderate: "

# Row 4 (Connected device): "Good" -> "High"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Characters(1, 4).Text = "High"

# Row 5 (Interview): "Good: " -> "High: ", "Medium: " -> "Moderate: "
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Characters(1, 6).Text = "High: "
$tbl.Cell(5, 4).Shape.TextFrame.TextRange.Characters(1, 8).Text = "Moderate: "
